$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

function Set-TextValue {
    param($cellRef, $text)
    $escaped = $text -replace '"', '""'
    $scratch.Formula = "=""" + $escaped + """"
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# Row 2
Set-TextValue "A2" "171670"
Set-TextValue "B2" "Монтаж трос стен вi-l дл.5986 серьга вверху (171670)"
Set-TextValue "C2" "111281"
Set-TextValue "D2" "165160, 161564"

# Row 3
Set-TextValue "A3" "111281"
Set-TextValue "B3" "Монтаж балка несущая вi-l5м дл. (111281)"
$ws.Range("C3").ClearContents()
Set-TextValue "D3" "171670, 192057"

# Row 4
Set-TextValue "A4" "100203"
Set-TextValue "B4" "Монтаж планка натягивающая дл.4990 100х40 (100203)"
Set-TextValue "C4" "165160, 161564"
$ws.Range("D4").ClearContents()

# Row 5
Set-TextValue "A5" "161564"
Set-TextValue "B5" "Монтаж панель стеклянная 1158х2735 мм с/с/а (161564)"
Set-TextValue "C5" "171670"
Set-TextValue "D5" "100203"

# Row 6
Set-TextValue "C6" "171670"

# Row 7
Set-TextValue "A7" "192057"
Set-TextValue "B7" "Монтаж алюминиевая направляющая для пола bi-level ii-уровень длина 4840 мм пр.285 к боковым стойкам (192057)"
Set-TextValue "C7" "111281"
$ws.Range("D7").ClearContents()

$scratch.ClearContents()
